{"js": "// Update the worksheet date and all 25 three-digit-by-one-digit\n// multiplication prompts to the new day's values.\n// Each entry is a unique, exact \"old text\" -> \"new text\" pair pulled\n// straight from the target diff, applied via Body.search()+insertText so\n// all run formatting (font/size) is preserved.\nconst replacements = [\n  [\"2024-05-22 Wednesday\", \"2024-05-23 Thursday\"],\n  [\"820\u00d78=\", \"348\u00d76=\"],\n  [\"227\u00d77=\", \"966\u00d72=\"],\n  [\"193\u00d76=\", \"456\u00d79=\"],\n  [\"474\u00d79=\", \"209\u00d78=\"],\n  [\"612\u00d74=\", \"333\u00d76=\"],\n  [\"655\u00d72=\", \"559\u00d79=\"],\n  [\"416\u00d75=\", \"327\u00d73=\"],\n  [\"731\u00d77=\", \"511\u00d78=\"],\n  [\"631\u00d72=\", \"376\u00d79=\"],\n  [\"537\u00d77=\", \"871\u00d79=\"],\n  [\"631\u00d78=\", \"272\u00d76=\"],\n  [\"362\u00d79=\", \"983\u00d74=\"],\n  [\"497\u00d74=\", \"843\u00d75=\"],\n  [\"132\u00d77=\", \"497\u00d77=\"],\n  [\"498\u00d73=\", \"689\u00d76=\"],\n  [\"121\u00d78=\", \"962\u00d78=\"],\n  [\"543\u00d79=\", \"489\u00d73=\"],\n  [\"132\u00d73=\", \"287\u00d78=\"],\n  [\"555\u00d78=\", \"647\u00d76=\"],\n  [\"824\u00d72=\", \"885\u00d78=\"],\n  [\"572\u00d79=\", \"837\u00d74=\"],\n  [\"809\u00d72=\", \"413\u00d79=\"],\n  [\"612\u00d73=\", \"926\u00d79=\"],\n  [\"194\u00d75=\", \"178\u00d76=\"],\n  [\"936\u00d76=\", \"528\u00d74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all 25 three-digit-by-one-digit\n# multiplication prompts to the new day's values using Find/Replace\n# (Word COM object model) so existing run formatting (font/size) on each\n# <w:t> is preserved.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-05-22 Wednesday\", \"2024-05-23 Thursday\"),\n    @(\"820\u00d78=\", \"348\u00d76=\"),\n    @(\"227\u00d77=\", \"966\u00d72=\"),\n    @(\"193\u00d76=\", \"456\u00d79=\"),\n    @(\"474\u00d79=\", \"209\u00d78=\"),\n    @(\"612\u00d74=\", \"333\u00d76=\"),\n    @(\"655\u00d72=\", \"559\u00d79=\"),\n    @(\"416\u00d75=\", \"327\u00d73=\"),\n    @(\"731\u00d77=\", \"511\u00d78=\"),\n    @(\"631\u00d72=\", \"376\u00d79=\"),\n    @(\"537\u00d77=\", \"871\u00d79=\"),\n    @(\"631\u00d78=\", \"272\u00d76=\"),\n    @(\"362\u00d79=\", \"983\u00d74=\"),\n    @(\"497\u00d74=\", \"843\u00d75=\"),\n    @(\"132\u00d77=\", \"497\u00d77=\"),\n    @(\"498\u00d73=\", \"689\u00d76=\"),\n    @(\"121\u00d78=\", \"962\u00d78=\"),\n    @(\"543\u00d79=\", \"489\u00d73=\"),\n    @(\"132\u00d73=\", \"287\u00d78=\"),\n    @(\"555\u00d78=\", \"647\u00d76=\"),\n    @(\"824\u00d72=\", \"885\u00d78=\"),\n    @(\"572\u00d79=\", \"837\u00d74=\"),\n    @(\"809\u00d72=\", \"413\u00d79=\"),\n    @(\"612\u00d73=\", \"926\u00d79=\"),\n    @(\"194\u00d75=\", \"178\u00d76=\"),\n    @(\"936\u00d76=\", \"528\u00d74=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap(=wdFindContinue), Format, ReplaceWith,\n    # Replace(=wdReplaceAll)\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
